$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1425304360311941
$ws.Range("D2").Value = 0.246097157014546
$ws.Range("E2").Value = 0.1773659978598765
$ws.Range("F2").Value = 1.012249720703693
$ws.Range("G2").Value = 0.5272962014676921
$ws.Range("H2").Value = 0.5812333897015094
$ws.Range("J2").Value = 0.1772330055721199
$ws.Range("N2").Value = 1.995241137946778
$ws.Range("O2").Value = 2.183072434228677
$ws.Range("B3").Value = 0.1329091933101978
$ws.Range("D3").Value = 0.2444676174134202
$ws.Range("E3").Value = 0.1742912681868667
$ws.Range("F3").Value = 0.9831260700589439
$ws.Range("G3").Value = 0.4999711970626919
$ws.Range("H3").Value = 0.5723711798586777
$ws.Range("J3").Value = 0.1721769985554289
$ws.Range("N3").Value = 1.846526605564804
$ws.Range("O3").Value = 2.106238768013412
$ws.Range("B4").Value = 0.1270661492628875
$ws.Range("D4").Value = 0.2435737801817908
$ws.Range("E4").Value = 0.1724976178116897
$ws.Range("F4").Value = 0.9658221175057946
$ws.Range("G4").Value = 0.4835087855673947
$ws.Range("H4").Value = 0.5672287518405454
$ws.Range("J4").Value = 0.1691764068403643
$ws.Range("N4").Value = 1.755520437593759
$ws.Range("O4").Value = 2.060343181411326
$ws.Range("B5").Value = 0.1247015078381395
$ws.Range("D5").Value = 0.2432364075533258
$ws.Range("E5").Value = 0.1717903946035122
$ws.Range("F5").Value = 0.9589159623540269
$ws.Range("G5").Value = 0.4768793569706844
$ws.Range("H5").Value = 0.5652084079914204
$ws.Range("J5").Value = 0.167979719279117
$ws.Range("N5").Value = 1.718516884911963
$ws.Range("O5").Value = 2.041962182227564
$ws.Range("B6").Value = 0.1243098614095572
$ws.Range("D6").Value = 0.24318201137892
$ws.Range("E6").Value = 0.171674392540087
$ws.Range("F6").Value = 0.957777978631924
$ws.Range("G6").Value = 0.4757833212346156
$ws.Range("H6").Value = 0.5648774770931766
$ws.Range("J6").Value = 0.1677825844652077
$ws.Range("N6").Value = 1.712377599108464
$ws.Range("O6").Value = 2.038929463360034
$ws.Range("B7").Value = 0.1270341919673399
$ws.Range("D7").Value = 0.2435691213927882
$ws.Range("E7").Value = 0.1724879839689208
$ws.Range("F7").Value = 0.9657283902091791
$ws.Range("G7").Value = 0.4834190585289946
$ws.Range("H7").Value = 0.5672012000805466
$ws.Range("J7").Value = 0.1691601623408303
$ws.Range("N7").Value = 1.755021054648495
$ws.Range("O7").Value = 2.060093986032228
$ws.Range("B8").Value = 0.1391997998654944
$ws.Range("D8").Value = 0.2455131683178706
$ws.Range("E8").Value = 0.1762862718187961
$ws.Range("F8").Value = 1.002087880873745
$ws.Range("G8").Value = 0.517809018632434
$ws.Range("H8").Value = 0.5781156669493015
$ws.Range("J8").Value = 0.1754681284190411
$ws.Range("N8").Value = 1.943904251618477
$ws.Range("O8").Value = 2.156314253543428
$ws.Range("B9").Value = 0.1635584291339143
$ws.Range("D9").Value = 0.2501707657315677
$ws.Range("E9").Value = 0.1844829058625663
$ws.Range("F9").Value = 1.077981530534032
$ws.Range("G9").Value = 0.5877591867547096
$ws.Range("H9").Value = 0.6018910408742215
$ws.Range("J9").Value = 0.1886638921441488
$ws.Range("N9").Value = 2.316504389646809
$ws.Range("O9").Value = 2.355184202377018
$ws.Range("B10").Value = 0.1817500012197257
$ws.Range("D10").Value = 0.2541069149559831
$ws.Range("E10").Value = 0.1909625588303214
$ws.Range("F10").Value = 1.136556074637113
$ws.Range("G10").Value = 0.6407026906574345
$ws.Range("H10").Value = 0.6208072618303788
$ws.Range("J10").Value = 0.1988664766722508
$ws.Range("N10").Value = 2.591327722520418
$ws.Range("O10").Value = 2.507549276033956
$ws.Range("B11").Value = 0.1900877858899008
$ws.Range("D11").Value = 0.2560090365128502
$ws.Range("E11").Value = 0.194010052836056
$ws.Range("F11").Value = 1.1638183204999
$ws.Range("G11").Value = 0.6651295264913699
$ws.Range("H11").Value = 0.6297279495998964
$ws.Range("J11").Value = 0.2036191076019378
$ws.Range("N11").Value = 2.716530669508643
$ws.Range("O11").Value = 2.578233826267081
$ws.Range("B12").Value = 0.1932538308938092
$ws.Range("D12").Value = 0.256745331136699
$ws.Range("E12").Value = 0.195178431529726
$ws.Range("F12").Value = 1.174230611703237
$ws.Range("G12").Value = 0.6744288402402958
$ws.Range("H12").Value = 0.6331513654079117
$ws.Range("J12").Value = 0.2054348799560159
$ws.Range("N12").Value = 2.763963193734583
$ws.Range("O12").Value = 2.605198180176444
$ws.Range("B13").Value = 0.1925715833279469
$ws.Range("D13").Value = 0.2565860458044114
$ws.Range("E13").Value = 0.1949261616900486
$ws.Range("F13").Value = 1.17198419430666
$ws.Range("G13").Value = 0.6724238672803722
$ws.Range("H13").Value = 0.6324120554357648
$ws.Range("J13").Value = 0.2050431063757969
$ws.Range("N13").Value = 2.753746911385861
$ws.Range("O13").Value = 2.599382130016579
$ws.Range("B14").Value = 0.190348085298524
$ws.Range("D14").Value = 0.2560692914184273
$ws.Range("E14").Value = 0.1941058882925901
$ws.Range("F14").Value = 1.164673168681375
$ws.Range("G14").Value = 0.6658935957561596
$ws.Range("H14").Value = 0.6300086875166357
$ws.Range("J14").Value = 0.2037681703561276
$ws.Range("N14").Value = 2.720432585385993
$ws.Range("O14").Value = 2.580448237419887
$ws.Range("B15").Value = 0.1889872547792919
$ws.Range("D15").Value = 0.2557548473068891
$ws.Range("E15").Value = 0.1936053170411043
$ws.Range("F15").Value = 1.160206502899541
$ws.Range("G15").Value = 0.6619000535299335
$ws.Range("H15").Value = 0.6285424593906441
$ws.Range("J15").Value = 0.202989326893956
$ws.Range("N15").Value = 2.700029152110062
$ws.Range("O15").Value = 2.568876427918497
$ws.Range("B16").Value = 0.1812063398625838
$ws.Range("D16").Value = 0.253984848216092
$ws.Range("E16").Value = 0.19076540734887
$ws.Range("F16").Value = 1.13478682902678
$ws.Range("G16").Value = 0.6391132448457597
$ws.Range("H16").Value = 0.6202306213541249
$ws.Range("J16").Value = 0.1985581258952323
$ws.Range("N16").Value = 2.583148644061055
$ws.Range("O16").Value = 2.502957520220775
$ws.Range("B17").Value = 0.1764487851809662
$ws.Range("D17").Value = 0.252927554079335
$ws.Range("E17").Value = 0.1890487917633763
$ws.Range("F17").Value = 1.119350579213645
$ws.Range("G17").Value = 0.6252221158806037
$ws.Range("H17").Value = 0.6152123840539616
$ws.Range("J17").Value = 0.1958682884033323
$ws.Range("N17").Value = 2.511489453100921
$ws.Range("O17").Value = 2.462870201680289
$ws.Range("B18").Value = 0.1737182494253204
$ws.Range("D18").Value = 0.2523299281001528
$ws.Range("E18").Value = 0.1880708399873683
$ws.Range("F18").Value = 1.110530072250853
$ws.Range("G18").Value = 0.6172645418874652
$ws.Range("H18").Value = 0.6123557350702526
$ws.Range("J18").Value = 0.1943316508989312
$ws.Range("N18").Value = 2.470290675910178
$ws.Range("O18").Value = 2.439942247904867
$ws.Range("B19").Value = 0.1727947553761737
$ws.Range("D19").Value = 0.2521293870371721
$ws.Range("E19").Value = 0.1877413368566323
$ws.Range("F19").Value = 1.107553563374182
$ws.Range("G19").Value = 0.614575776360283
$ws.Range("H19").Value = 0.6113936264419806
$ws.Range("J19").Value = 0.19381317184191
$ws.Range("N19").Value = 2.456344689235891
$ws.Range("O19").Value = 2.432201422091509
$ws.Range("B20").Value = 0.1769546279149097
$ws.Range("D20").Value = 0.2530390182400026
$ws.Range("E20").Value = 0.1892305555620766
$ws.Range("F20").Value = 1.120987789708067
$ws.Range("G20").Value = 0.6266975129551611
$ws.Range("H20").Value = 0.6157435095531412
$ws.Range("J20").Value = 0.1961535407199193
$ws.Range("N20").Value = 2.519115905986666
$ws.Range("O20").Value = 2.467124190622599
$ws.Range("B21").Value = 0.1910009463247206
$ws.Range("D21").Value = 0.256220640706772
$ws.Range("E21").Value = 0.1943464327933526
$ws.Range("F21").Value = 1.166818187367696
$ws.Range("G21").Value = 0.6678103529505393
$ws.Range("H21").Value = 0.6307133845868691
$ws.Range("J21").Value = 0.2041422140131317
$ws.Range("N21").Value = 2.730217288371648
$ws.Range("O21").Value = 2.586004214704531
$ws.Range("B22").Value = 0.2002316580428385
$ws.Range("D22").Value = 0.2583932720500997
$ws.Range("E22").Value = 0.1977736478161205
$ws.Range("F22").Value = 1.197287856512247
$ws.Range("G22").Value = 0.6949679831682829
$ws.Range("H22").Value = 0.6407613510679369
$ws.Range("J22").Value = 0.2094568644435526
$ws.Range("N22").Value = 2.868303394937357
$ws.Range("O22").Value = 2.664851336729612
$ws.Range("B23").Value = 0.1953005061487687
$ws.Range("D23").Value = 0.2572251772788974
$ws.Range("E23").Value = 0.1959368213278907
$ws.Range("F23").Value = 1.180978317030494
$ws.Range("G23").Value = 0.6804470454142972
$ws.Range("H23").Value = 0.6353743924378534
$ws.Range("J23").Value = 0.206611761512363
$ws.Range("N23").Value = 2.794595190495954
$ws.Range("O23").Value = 2.622663646788794
$ws.Range("B24").Value = 0.1767259219246853
$ws.Range("D24").Value = 0.2529885934366405
$ws.Range("E24").Value = 0.1891483522605455
$ws.Range("F24").Value = 1.120247438639538
$ws.Range("G24").Value = 0.6260303967945617
$ws.Range("H24").Value = 0.6155032992347458
$ws.Range("J24").Value = 0.1960245476698788
$ws.Range("N24").Value = 2.51566798917645
$ws.Range("O24").Value = 2.465200592235362
$ws.Range("B25").Value = 0.1569161479396968
$ws.Range("D25").Value = 0.2488203805644389
$ws.Range("E25").Value = 0.1821852480318853
$ws.Range("F25").Value = 1.056956987943394
$ws.Range("G25").Value = 0.5685648655344551
$ws.Range("H25").Value = 0.5952049802152146
$ws.Range("J25").Value = 0.1850052413101935
$ws.Range("N25").Value = 2.215498394812073
$ws.Range("O25").Value = 2.300289579330752
